$d = $word.ActiveDocument

# Helper constants for Find/Execute parameter positions:
# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# wdReplace: 0 = wdReplaceNone, 1 = wdReplaceOne, 2 = wdReplaceAll
# wdFindWrap: 1 = wdFindContinue

# --- Bullet 1: "Not much work done..." split and expanded ---
$d.Content.Find.Execute(
    "Not much work done, first formal meeting, first informal meeting. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "First formal meeting, first informal meeting. Plans for project plan discussed and jobs allocated. A lot of time was spent discussing hour allocation in meetings.",
    2) | Out-Null

# --- Bullet 2: "Deadline for project plan..." gains (D1) and deadlines/are ---
$d.Content.Find.Execute(
    "Deadline for project plan.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deadline for project plan (D1).",
    2) | Out-Null

$d.Content.Find.Execute(
    "making sure deadline is being stuck to",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "making sure deadlines are being stuck to",
    2) | Out-Null

# --- Bullet 3: "Establish a proper folder..." -> "Established..." + working pair names + UML sentence ---
$d.Content.Find.Execute(
    "Establish a proper folder structure for SVN. Working pairs was appointed. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Established a proper folder structure for SVN. Working pairs was appointed: Kelvin and Oscar on Model, Kristian and Edward on View and Brian on Controller. ",
    2) | Out-Null

$d.Content.Find.Execute(
    "was discussed in formal and informal meeting.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "was discussed in formal and informal meeting. UML work begun.",
    2) | Out-Null

# --- Bullet 4: "Coding begins on Model and View..." expanded ---
$d.Content.Find.Execute(
    "Coding begins on Model and View. Formal meeting reviewed progress. Controller person helps out with View.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Coding begins on Model and View by the working pairs. Formal meeting reviewed progress. Brian the Controller person helps out with View offering design consultancy as View is highly tied with Controller.",
    2) | Out-Null

# --- Bullet 5: "Discussed ideas for D2..." gains UML diagram sentence ---
$d.Content.Find.Execute(
    "Discussed ideas for D2. Prepared code and presentation for D2.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discussed ideas for D2. Prepared code and presentation for D2. UML diagrams were finalised. A long library group meeting was held on the Sunday afternoon for this work to be done as a group because the presentation will be presented by the whole group and needs to be dynamic.",
    2) | Out-Null

# --- Bullet 6: "Code freeze. Rehearsal..." -> "Final rehearsal..." + extra sentence ---
$d.Content.Find.Execute(
    "Code freeze. Rehearsal and presentation of D2.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Code freeze. Final rehearsal and presentation of D2. Final rehearsal was done in an informal meeting an hour before the presentation so everyone knew what they are presenting and were fresh with the presentation in mind before the presentation.",
    2) | Out-Null

# --- Bullet 7: "Reviewed coding progress..." gains extra sentence ---
$d.Content.Find.Execute(
    "Reviewed coding progress. Work begins on D3.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Reviewed coding progress. Work begins on D3. D2 marks were received and met with approving nods all round from members of the group. Controller code begins. Confidence level was high.",
    2) | Out-Null

# --- Bullet 8: "Easter vacation..." wrapped with ** ... ** ---
$d.Content.Find.Execute(
    "Easter vacation: no work was done.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "** Easter vacation: no work was done. **",
    2) | Out-Null

# --- Bullet 9: "Splitting /src folder..." tense change + "more" ---
$d.Content.Find.Execute(
    "folder into model and view is no longer a good solution",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "folder into model and view was no longer a good solution",
    2) | Out-Null

$d.Content.Find.Execute(
    "starting to be developed. /",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "starting to be developed more. /",
    2) | Out-Null

# --- Bullet "Not establishing interfaces..." : people -> groups ---
$d.Content.Find.Execute(
    "both Model and View people to get an idea",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "both Model and View groups to get an idea",
    2) | Out-Null

# --- Restore the "_GoBack" bookmark at its new location (moved by the edit),
#     sitting between "derail " and "progress." in the final "Lesson learnt" bullet ---
$gb = $d.Content
$gb.Find.Execute("derail ") | Out-Null
$gbRange = $d.Range($gb.End, $gb.End)
$d.Bookmarks.Add("_GoBack", $gbRange) | Out-Null
